# Applies the crypto price/volume refresh described in the commit:
# "Updated cryptos list on Sun May 19 09:18:03 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 (Bitcoin)
$ws.Range("D2").Value = '67.095.48'
$ws.Range("E2").Value = '  +0.07%  '

# Row 3 (Ethereum)
$ws.Range("D3").Value = '3.105.97'
$ws.Range("E3").Value = '  -0.34%  '

# Row 4 (TetherUSD)
$ws.Range("E4").Value = '  +0.08%  '

# Row 5 (BNB)
$ws.Range("D5").Value = '''579.48'
$ws.Range("E5").Value = '  -0.10%  '

# Row 6 (Solana)
$ws.Range("D6").Value = '''172.62'
$ws.Range("E6").Value = '  -0.43%  '

# Row 7 (USDC)
$ws.Range("E7").Value = '  +0.05%  '

# Row 8 (XRP)
$ws.Range("E8").Value = '  -0.87%  '

# Row 9 (Toncoin)
$ws.Range("D9").Value = '''6.51'
$ws.Range("E9").Value = '  +1.12%  '

# Row 10 (Dogecoin)
$ws.Range("E10").Value = '  -1.68%  '

# Row 11 (Cardano)
$ws.Range("D11").Value = '''0.476'
$ws.Range("E11").Value = '  -1.27%  '

# Row 12 (ShibaInu)
$ws.Range("D12").Value = '''0.0000246'
$ws.Range("E12").Value = '  -1.31%  '

# Row 13 (Avalanche)
$ws.Range("D13").Value = '''36.64'
$ws.Range("E13").Value = '  -1.60%  '

# Row 14 (TRON)
$ws.Range("E14").Value = '  -1.69%  '

# Row 15 (WrappedliquidstakedEther2.0)
$ws.Range("D15").Value = '3.623.31'
$ws.Range("E15").Value = '  -0.23%  '

# Row 16 (WrappedBTC)
$ws.Range("D16").Value = '67.091.84'
$ws.Range("E16").Value = '  +0.13%  '

# Row 17 (Polkadot)
$ws.Range("D17").Value = '''7.07'
$ws.Range("E17").Value = '  -1.69%  '

# Row 18 (WrappedEther)
$ws.Range("D18").Value = '3.106.57'
$ws.Range("E18").Value = '  -0.27%  '

# Row 19 (Chainlink)
$ws.Range("D19").Value = '''16.64'
$ws.Range("E19").Value = '  +2.73%  '

# Row 20 (BitcoinCash)
$ws.Range("D20").Value = '''490.55'
$ws.Range("E20").Value = '  +0.66%  '

# Row 21 (Polygon)
$ws.Range("D21").Value = '''0.703'
$ws.Range("E21").Value = '  -2.17%  '

# Row 22 (Uniswap)
$ws.Range("D22").Value = '''7.83'
$ws.Range("E22").Value = '  +2.79%  '

# Row 23 (Litecoin)
$ws.Range("D23").Value = '''83.87'
$ws.Range("E23").Value = '  -0.72%  '

# Row 24 (InternetComputer(DFINITY))
$ws.Range("D24").Value = '''13.05'
$ws.Range("E24").Value = '  -2.40%  '

# Row 25 (Fetch.AI)
$ws.Range("E25").Value = '  -3.59%  '

# Row 26 (RenderToken)
$ws.Range("D26").Value = '''10.56'
$ws.Range("E26").Value = '  +4.58%  '

# Row 27 (Dai)
$ws.Range("E27").Value = '  -0.03%  '

# Row 28 (NEARProtocol)
$ws.Range("D28").Value = '''7.88'
$ws.Range("E28").Value = '  -2.02%  '

# Row 29 (ImmutableX)
$ws.Range("E29").Value = '  -2.88%  '

# Row 30 (PancakeSwap)
$ws.Range("E30").Value = '  -1.23%  '

# Row 31 (EthereumClassic)
$ws.Range("D31").Value = '''28.23'
$ws.Range("E31").Value = '  -2.39%  '

# Row 32 (Hedera)
$ws.Range("D32").Value = '''0.113'

# Row 33 (PEPE)
$ws.Range("D33").Value = '0.0₃0925'
$ws.Range("E33").Value = '  -7.65%  '

# Row 34 (FirstDigitalUSD)
$ws.Range("D34").Value = '''1.00'
$ws.Range("E34").Value = '  +0.10%  '

# Row 35 (Filecoin)
$ws.Range("D35").Value = '''5.81'
$ws.Range("E35").Value = '  -2.21%  '

# Row 36 (Mantle)
$ws.Range("D36").Value = '''0.973'
$ws.Range("E36").Value = '  -1.83%  '

# Row 37 (Arweave)
$ws.Range("D37").Value = '''46.93'
$ws.Range("E37").Value = '  -1.24%  '

# Row 38 (Stacks)
$ws.Range("D38").Value = '''2.03'
$ws.Range("E38").Value = '  -4.03%  '

# Row 39 (TheGraph)
$ws.Range("B39").Value = 'Kaspa'
$ws.Range("C39").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D39").Value = '''0.124'
$ws.Range("E39").Value = '  +0.69%  '

# Row 40 (Kaspa)
$ws.Range("B40").Value = 'TheGraph'
$ws.Range("C40").Value = 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt'
$ws.Range("D40").Value = '''0.307'
$ws.Range("E40").Value = '  -3.02%  '

# Row 41 (Cosmos)
$ws.Range("D41").Value = '''8.44'
$ws.Range("E41").Value = '  -2.62%  '

# Row 42 (Bittensor)
$ws.Range("D42").Value = '''386.23'
$ws.Range("E42").Value = '  +0.32%  '

# Row 43 (Maker)
$ws.Range("D43").Value = '2.802.06'
$ws.Range("E43").Value = '  -1.53%  '

# Row 44 (dogwifhat)
$ws.Range("E44").Value = '  -8.94%  '

# Row 45 (VeChain)
$ws.Range("E45").Value = '  -2.47%  '

# Row 46 (Monero)
$ws.Range("D46").Value = '''134.84'
$ws.Range("E46").Value = '  -1.65%  '

# Row 47 (USDe)
$ws.Range("E47").Value = '  -0.03%  '

# Row 48 (InjectiveProtocol)
$ws.Range("D48").Value = '''24.94'
$ws.Range("E48").Value = '  -0.79%  '

# Row 49 (ThetaToken)
$ws.Range("D49").Value = '''2.19'
$ws.Range("E49").Value = '  -1.98%  '

# Row 50 (Stellar)
$ws.Range("E50").Value = '  -1.75%  '

# Row 51 (THORChain)
$ws.Range("D51").Value = '''6.69'
$ws.Range("E51").Value = '  -2.47%  '
